# The edit inserts one new data row (a new Cereza "Rainier" price record for
# Región de O'Higgins, dated 2023-11-28 / serial 45258) right before the
# existing row 322, pushing all subsequent rows (old 322-400) down by one
# (new rows 323-401). The workbook dimension grows from A1:T400 to A1:T401.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 322, shifting row 322 and everything below
# it down by one row.
$ws.Rows(322).Insert()

# Populate the newly inserted row 322 with the new record. Columns A, B, C,
# E, F, G, H, I, J hold the same constant market/product metadata used by
# every other row in this sheet.
$ws.Range("A322").Value = 10
$ws.Range("B322").Value = "Vega Modelo de Temuco"
$ws.Range("C322").Value = "La Araucanía"
$ws.Range("D322").Value = 45258
$ws.Range("E322").Value = 9
$ws.Range("F322").Value = "Fruta"
$ws.Range("G322").Value = 100103
$ws.Range("H322").Value = "Frutos de hueso (carozo)"
$ws.Range("I322").Value = 100103001
$ws.Range("J322").Value = "Cereza"
$ws.Range("K322").Value = "Rainier"
$ws.Range("L322").Value = "Primera"
$ws.Range("M322").Value = 350
$ws.Range("N322").Value = 25000
$ws.Range("O322").Value = 25000
$ws.Range("P322").Value = 25000
$ws.Range("Q322").Value = "`$/bandeja 10 kilos"
$ws.Range("R322").Value = "Región de O'Higgins"
$ws.Range("S322").Value = 2500
$ws.Range("T322").Value = 10
